# Update column G ("K") values on Sheet1 to reflect the regenerated
# save_data (K replaces the old "Strike#" metric; std/mean recalculated
# and s_vals rewritten upstream). Only the numeric values in column G
# change for data rows 2-40 (row 13 keeps its original value of 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 1
    6  = 1
    7  = 2
    8  = 2
    9  = 3
    10 = 2
    11 = 0
    12 = 3
    14 = 3
    15 = 2
    16 = 0
    17 = 3
    18 = 2
    19 = 2
    20 = 1
    21 = 1
    22 = 0
    23 = 1
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 0
    29 = 3
    30 = 0
    31 = 1
    32 = 3
    33 = 2
    34 = 2
    35 = 3
    36 = 1
    37 = 2
    38 = 1
    39 = 0
    40 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
